$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 525
$ws.Range("F4").Value = 559
$ws.Range("F5").Value = 9203
$ws.Range("F7").Value = 11908
$ws.Range("G7").Value = 238
$ws.Range("F8").Value = 11908
$ws.Range("G8").Value = 238
$ws.Range("F9").Value = 105
$ws.Range("F14").Value = 37
$ws.Range("F15").Value = 120
$ws.Range("F17").Value = 442
$ws.Range("F18").Value = 2031
$ws.Range("F19").Value = 801
$ws.Range("F20").Value = 760
$ws.Range("F23").Value = 403
$ws.Range("F25").Value = 92
$ws.Range("F27").Value = 9
$ws.Range("F28").Value = 1487
$ws.Range("F30").Value = 18
$ws.Range("F31").Value = 15
$ws.Range("F32").Value = 54
$ws.Range("F33").Value = 503
$ws.Range("F36").Value = 498
$ws.Range("F37").Value = 335
$ws.Range("F38").Value = 485
$ws.Range("F39").Value = 377
$ws.Range("F40").Value = 2117
$ws.Range("G40").Value = 72
$ws.Range("F41").Value = 72
$ws.Range("F43").Value = 564
$ws.Range("F44").Value = 437
$ws.Range("F45").Value = 151
$ws.Range("F46").Value = 872
$ws.Range("F49").Value = 264
$ws.Range("F50").Value = 239

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 42
$ws.Range("F7").Value = 73
$ws.Range("F13").Value = 45
$ws.Range("F24").Value = 107
$ws.Range("F25").Value = 63
$ws.Range("F26").Value = 424

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2902
$ws.Range("F4").Value = 364
$ws.Range("F5").Value = 234
$ws.Range("F6").Value = 236

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 525
$ws.Range("F4").Value = 42
$ws.Range("F5").Value = 364
$ws.Range("F6").Value = 234
$ws.Range("F7").Value = 559
$ws.Range("F8").Value = 9203
$ws.Range("F10").Value = 11908
$ws.Range("G10").Value = 238
$ws.Range("F13").Value = 37
$ws.Range("F14").Value = 120
$ws.Range("F16").Value = 2031
$ws.Range("F17").Value = 801
$ws.Range("F18").Value = 760
$ws.Range("F21").Value = 403
$ws.Range("F25").Value = 73
$ws.Range("F27").Value = 9
$ws.Range("F28").Value = 236
$ws.Range("F29").Value = 1487
$ws.Range("F31").Value = 45
$ws.Range("F32").Value = 503
$ws.Range("F37").Value = 498
$ws.Range("F38").Value = 485
$ws.Range("F39").Value = 377
$ws.Range("F40").Value = 2117
$ws.Range("G40").Value = 72
$ws.Range("C41").Value = "北京·开饭咯！迷宫饭同人ONLY"
$ws.Range("D41").Value = "酒仙桥北路2号院798艺术区706后街1号 北京格瑞斯艺术酒店"
$ws.Range("E41").Value = "2024.11.09 10:30-11.10 16:30"
$ws.Range("F41").Value = 72
$ws.Range("G41").Value = 19.9
$ws.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=92529"
$ws.Range("I41").Value = "//i2.hdslb.com/bfs/openplatform/202409/Hq5V5Geo1727434065793.png"
$ws.Range("F42").Value = 564
$ws.Range("F43").Value = 438
$ws.Range("F44").Value = 151
$ws.Range("F46").Value = 63
$ws.Range("F47").Value = 424
$ws.Range("F49").Value = 265
$ws.Range("F50").Value = 239
